$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing participant ID (B2): p0410947 -> p04100947 ---
$ws.Range("B2").Value = "p04100947"

# --- Row 18 ---
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("O18").PasteSpecial(-4122)
$ws.Range("Q18").PasteSpecial(-4122)
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("S18").PasteSpecial(-4122)
$ws.Range("T18").PasteSpecial(-4122)
$ws.Range("U18").PasteSpecial(-4122)
$ws.Range("V18").PasteSpecial(-4122)
$ws.Range("Y18").PasteSpecial(-4122)
$ws.Range("A18").Value = 45398.458225312497
$ws.Range("B18").Value = "p04161322"
$ws.Range("C18").Value = "Yes"
$ws.Range("D18").Value = "No"
$ws.Range("F18").Value = "4 - 5 years of programming experience"
$ws.Range("G18").Value = "No experience"
$ws.Range("H18").Value = "No experience"
$ws.Range("I18").Value = "No experience"
$ws.Range("J18").Value = "Junior"
$ws.Range("K18").Value = "Lecture Materials (e.g., In-Class Slides, Lecture Notes)"
$ws.Range("L18").Value = "The materials do not provide the sufficient or required information"
$ws.Range("M18").Value = "Lecture Materials (e.g., In-Class Slides, Lecture Notes)"
$ws.Range("N18").Value = "The materials do not provide the sufficient or required information, The information and resources are not always accessible"
$ws.Range("O18").Value = "Chat"
$ws.Range("Q18").Value = "Neither Agree or Disagree"
$ws.Range("R18").Value = "Agree"
$ws.Range("S18").Value = "Agree"
$ws.Range("T18").Value = "Agree"
$ws.Range("U18").Value = "Disagree"
$ws.Range("V18").Value = "Agree"

# --- Row 19 ---
$ws.Range("A17").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("O19").PasteSpecial(-4122)
$ws.Range("P19").PasteSpecial(-4122)
$ws.Range("Q19").PasteSpecial(-4122)
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("T19").PasteSpecial(-4122)
$ws.Range("U19").PasteSpecial(-4122)
$ws.Range("V19").PasteSpecial(-4122)
$ws.Range("X19").PasteSpecial(-4122)
$ws.Range("Y19").PasteSpecial(-4122)
$ws.Range("Z19").PasteSpecial(-4122)
$ws.Range("A19").Value = 45398.637042881943
$ws.Range("B19").Value = "p04161834"
$ws.Range("C19").Value = "Yes"
$ws.Range("D19").Value = "No"
$ws.Range("E19").Value = "Participant left early, spending 50 minutes in the experiment. He partially completed the task."
$ws.Range("F19").Value = "1 - 3 years of programming experience"
$ws.Range("G19").Value = "1 - 3 years of experience"
$ws.Range("H19").Value = "No experience"
$ws.Range("I19").Value = "No experience"
$ws.Range("J19").Value = "Senior"
$ws.Range("K19").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Lecture Materials (e.g., In-Class Slides, Lecture Notes), Chatbots (e.g., ChatGPT), Technical Documentation"
$ws.Range("L19").Value = "The materials do not provide the sufficient or required information, The information and resources are not in a form that is readily useable, The information and resources are not concise and clear, The information and resources are not organized into logical and understandable components"
$ws.Range("M19").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Lecture Materials (e.g., In-Class Slides, Lecture Notes), Chatbots (e.g., ChatGPT), Technical Documentation"
$ws.Range("N19").Value = "The information and resources are not concise and clear, The information and resources are not organized into logical and understandable components"
$ws.Range("O19").Value = "I did not use the desktop computer"
$ws.Range("P19").Value = "I did not use the desktop computer"
$ws.Range("Q19").Value = "I did not use the desktop computer"
$ws.Range("R19").Value = "I did not use the desktop computer"
$ws.Range("S19").Value = "I did not use the desktop computer"
$ws.Range("T19").Value = "I did not use the desktop computer"
$ws.Range("U19").Value = "I did not use the desktop computer"
$ws.Range("V19").Value = "I did not use the desktop computer"
$ws.Range("X19").Value = "I think my experience was good, I may have liked to have a partner or a 2nd thought on things."
$ws.Range("Z19").Value = "I had a positive experience with the time I had, I had to leave early but putting together the block based language to see the robot arm moving the coffee cans was pretty cool."

# --- Row 20 ---
$ws.Range("A17").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("J20").PasteSpecial(-4122)
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("L20").PasteSpecial(-4122)
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("O20").PasteSpecial(-4122)
$ws.Range("P20").PasteSpecial(-4122)
$ws.Range("Q20").PasteSpecial(-4122)
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("T20").PasteSpecial(-4122)
$ws.Range("U20").PasteSpecial(-4122)
$ws.Range("V20").PasteSpecial(-4122)
$ws.Range("X20").PasteSpecial(-4122)
$ws.Range("Z20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45399.515523819442
$ws.Range("B20").Value = "p04171512"
$ws.Range("C20").Value = "Yes"
$ws.Range("D20").Value = "No"
$ws.Range("F20").Value = "1 - 3 years of programming experience"
$ws.Range("G20").Value = "Less than 1 year of experience"
$ws.Range("H20").Value = "Less than 1 year of experience"
$ws.Range("I20").Value = "Limited experience"
$ws.Range("J20").Value = "Junior"
$ws.Range("K20").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Chatbots (e.g., ChatGPT)"
$ws.Range("L20").Value = "The information and resources are not in a form that is readily useable, The information and resources are not organized into logical and understandable components"
$ws.Range("M20").Value = "Videos (e.g., YouTube, Udemy), Chatbots (e.g., ChatGPT), Technical Documentation"
$ws.Range("N20").Value = "Needs more video representation"
$ws.Range("O20").Value = "Chat"
$ws.Range("P20").Value = "The chat bot was very responsive"
$ws.Range("Q20").Value = "Agree"
$ws.Range("R20").Value = "Agree"
$ws.Range("S20").Value = "Agree"
$ws.Range("T20").Value = "Neither Agree or Disagree"
$ws.Range("U20").Value = "Neither Agree or Disagree"
$ws.Range("V20").Value = "Agree"
$ws.Range("X20").Value = "Yes defintely couldve been clearer what the task is"
$ws.Range("Z20").Value = "The block language is pretty intuitive"

# --- Row 21 ---
$ws.Range("A17").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("I21").PasteSpecial(-4122)
$ws.Range("J21").PasteSpecial(-4122)
$ws.Range("K21").PasteSpecial(-4122)
$ws.Range("M21").PasteSpecial(-4122)
$ws.Range("O21").PasteSpecial(-4122)
$ws.Range("P21").PasteSpecial(-4122)
$ws.Range("Q21").PasteSpecial(-4122)
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("S21").PasteSpecial(-4122)
$ws.Range("T21").PasteSpecial(-4122)
$ws.Range("U21").PasteSpecial(-4122)
$ws.Range("V21").PasteSpecial(-4122)
$ws.Range("X21").PasteSpecial(-4122)
$ws.Range("Y21").PasteSpecial(-4122)
$ws.Range("Z21").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("A21").Value = 45399.66549944444
$ws.Range("B21").Value = "p04171840"
$ws.Range("C21").Value = "Yes"
$ws.Range("D21").Value = "No"
$ws.Range("E21").Value = "The proctor put the wrong ID for this participant (p04171950 is p04171840). The program would not let him apply any new work and froze for a bit. His code is based off his logic, not testing since the robot would not update with new code in the last 15 minutes."
$ws.Range("F21").Value = "1 - 3 years of programming experience"
$ws.Range("G21").Value = "No experience"
$ws.Range("H21").Value = "No experience"
$ws.Range("I21").Value = "No experience"
$ws.Range("J21").Value = "Senior"
$ws.Range("K21").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Chatbots (e.g., ChatGPT)"
$ws.Range("M21").Value = "Videos (e.g., YouTube, Udemy)"
$ws.Range("O21").Value = "Watch"
$ws.Range("P21").Value = "Videos were very helpful and programming the robot was easy."
$ws.Range("Q21").Value = "Strongly Agree"
$ws.Range("R21").Value = "Strongly Agree"
$ws.Range("S21").Value = "Strongly Agree"
$ws.Range("T21").Value = "Strongly Agree"
$ws.Range("U21").Value = "Strongly Agree"
$ws.Range("V21").Value = "Strongly Agree"
$ws.Range("X21").Value = "none"
$ws.Range("Z21").Value = "The language is fine, just the robot did not apply the new code I gave him"

# --- Row 22 ---
$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("O22").PasteSpecial(-4122)
$ws.Range("P22").PasteSpecial(-4122)
$ws.Range("Q22").PasteSpecial(-4122)
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("S22").PasteSpecial(-4122)
$ws.Range("T22").PasteSpecial(-4122)
$ws.Range("U22").PasteSpecial(-4122)
$ws.Range("V22").PasteSpecial(-4122)
$ws.Range("X22").PasteSpecial(-4122)
$ws.Range("Z22").PasteSpecial(-4122)
$ws.Range("A22").Value = 45400.500941493054
$ws.Range("B22").Value = "p04181446"
$ws.Range("C22").Value = "Yes"
$ws.Range("D22").Value = "Yes"
$ws.Range("F22").Value = "4 - 5 years of programming experience"
$ws.Range("G22").Value = "No experience"
$ws.Range("H22").Value = "No experience"
$ws.Range("I22").Value = "Limited experience"
$ws.Range("J22").Value = "PhD student"
$ws.Range("K22").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Lecture Materials (e.g., In-Class Slides, Lecture Notes), Chatbots (e.g., ChatGPT), Technical Documentation"
$ws.Range("L22").Value = "The materials do not provide the sufficient or required information, The information and resources are not always accessible, The information and resources are not in a form that is readily useable, The information and resources are not concise and clear, The information and resources are not organized into logical and understandable components, The information and resources provided are not up to date"
$ws.Range("M22").Value = "I did not study robot pogramming"
$ws.Range("N22").Value = "The information and resources are not organized into logical and understandable components"
$ws.Range("O22").Value = "I did not use the desktop computer"
$ws.Range("P22").Value = "I did not use them"
$ws.Range("Q22").Value = "I did not use the desktop computer"
$ws.Range("R22").Value = "I did not use the desktop computer"
$ws.Range("S22").Value = "I did not use the desktop computer"
$ws.Range("T22").Value = "I did not use the desktop computer"
$ws.Range("U22").Value = "I did not use the desktop computer"
$ws.Range("V22").Value = "I did not use the desktop computer"
$ws.Range("X22").Value = "The experiment was very well setup. I like that it required repetition while also changing some values that were repeated (i.e. location that the robot would put the cans in).
Additionally, I liked that the experiment was challenging but also accessible. Especially for someone with some coding experience, but no robot coding experience."
$ws.Range("Z22").Value = "I feel that the block based language is very easy to understand, especially for people who have a basic understanding of coding. 
The block based language makes creating functions much easier, as repetitive tasks can be copied and dragged over. However, a user interacting with the block based language may struggle to understand how to duplicate blocks without guidance. 
I think that the block based programming language would benefit a lot from a smoother display/machine that it is running on. Some issues I ran into would be maybe expecting to be able to pinch to zoom, or attempt to highlight a block and drag which I could not do. However this is a hardware issue, and not at all an issue with the block based language. "

# --- Row 23 ---
$ws.Range("A17").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("K23").PasteSpecial(-4122)
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("O23").PasteSpecial(-4122)
$ws.Range("Q23").PasteSpecial(-4122)
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("S23").PasteSpecial(-4122)
$ws.Range("T23").PasteSpecial(-4122)
$ws.Range("U23").PasteSpecial(-4122)
$ws.Range("V23").PasteSpecial(-4122)
$ws.Range("Y23").PasteSpecial(-4122)
$ws.Range("Z23").PasteSpecial(-4122)
$ws.Range("A23").Value = 45400.57352726852
$ws.Range("B23").Value = "p04181705"
$ws.Range("C23").Value = "Yes"
$ws.Range("D23").Value = "Yes"
$ws.Range("F23").Value = "More than 5 years of programming experience"
$ws.Range("G23").Value = "Less than 1 year of experience"
$ws.Range("H23").Value = "No experience"
$ws.Range("I23").Value = "Limited experience"
$ws.Range("J23").Value = "Senior"
$ws.Range("K23").Value = "Videos (e.g., YouTube, Udemy), Online Communities (e.g., Reddit, Stack Overflow), Lecture Materials (e.g., In-Class Slides, Lecture Notes), Chatbots (e.g., ChatGPT), Technical Documentation"
$ws.Range("L23").Value = "The information and resources are not in a form that is readily useable, The information and resources are not concise and clear, The information and resources provided are not up to date"
$ws.Range("O23").Value = "I did not use the desktop computer"
$ws.Range("Q23").Value = "I did not use the desktop computer"
$ws.Range("R23").Value = "I did not use the desktop computer"
$ws.Range("S23").Value = "I did not use the desktop computer"
$ws.Range("T23").Value = "I did not use the desktop computer"
$ws.Range("U23").Value = "I did not use the desktop computer"
$ws.Range("V23").Value = "I did not use the desktop computer"
$ws.Range("Z23").Value = "Even with my last experience with Scratch happened 7 years ago It was super easy to use the tool provided."

$excel.CutCopyMode = $false
